$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 59
$wsExhibit.Range("F3").Value = 65
$wsExhibit.Range("F4").Value = 158
$wsExhibit.Range("F5").Value = 354
$wsExhibit.Range("F6").Value = 5232
$wsExhibit.Range("F7").Value = 115
$wsExhibit.Range("F8").Value = 5334
$wsExhibit.Range("F9").Value = 616
$wsExhibit.Range("F10").Value = 3
$wsExhibit.Range("F11").Value = 1359

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 59
$wsAll.Range("F3").Value = 65
$wsAll.Range("F4").Value = 158
$wsAll.Range("F6").Value = 354
$wsAll.Range("F7").Value = 5232
$wsAll.Range("F8").Value = 115
$wsAll.Range("F9").Value = 5334
$wsAll.Range("F10").Value = 616
$wsAll.Range("F11").Value = 3
$wsAll.Range("F12").Value = 1359
